$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the underlying input values (row 19-23 of the "balance" table) ---
# Row 19
$ws.Range("E19").Value = 1.4

# Row 20
$ws.Range("E20").Value = 1.55

# Row 21
$ws.Range("D21").Value = 1.6
$ws.Range("E21").Value = 1.6
$ws.Range("F21").Value = 1.15

# Row 22
$ws.Range("D22").Value = 1.7
$ws.Range("E22").Value = 1.75
$ws.Range("F22").Value = 1.15

# Row 23
$ws.Range("B23").Value = 1.8
$ws.Range("C23").Value = 1.8
$ws.Range("F23").Value = 1.15

# --- Update the active selection to match the saved view state ---
$ws.Range("B24").Select() | Out-Null
